$d = $word.ActiveDocument
$sec = $d.Sections(1)

# --- Footer 1 (default footer) - Pearson logo: image2.png -> image1.png ---
$ftr1 = $sec.Footers(1)
$xml1 = $ftr1.Range.WordOpenXML
$xml1 = $xml1.Replace('id="1" name="image2.png"', 'id="1" name="image1.png"')
$xml1 = $xml1.Replace('id="0" name="image2.png"', 'id="0" name="image1.png"')
$ftr1.Range.WordOpenXML = $xml1

# --- Footer 2 (first-page footer) - Pearson logo: image2.png -> image1.png ---
$ftr2 = $sec.Footers(2)
$xml2 = $ftr2.Range.WordOpenXML
$xml2 = $xml2.Replace('id="2" name="image2.png"', 'id="2" name="image1.png"')
$xml2 = $xml2.Replace('id="0" name="image2.png"', 'id="0" name="image1.png"')
$ftr2.Range.WordOpenXML = $xml2

# --- Header 2 (first-page header) - BTec logo: image1.jpg -> image2.jpg ---
$hdr2 = $sec.Headers(2)
$xml3 = $hdr2.Range.WordOpenXML
$xml3 = $xml3.Replace('id="3" name="image1.jpg"', 'id="3" name="image2.jpg"')
$xml3 = $xml3.Replace('id="0" name="image1.jpg"', 'id="0" name="image2.jpg"')
$hdr2.Range.WordOpenXML = $xml3

Write-Host "done"
